$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H4").Value = 255.77777
$ws_ALC.Range("J4").Value = 400
$ws_ALC.Range("L4").Value = 400
$ws_ALC.Range("N4").Value = -628
$ws_ALC.Range("H86").Value = 34788.445
$ws_ALC.Range("I86").Value = 17757.666
$ws_ALC.Range("J86").Value = 68850
$ws_ALC.Range("K86").Value = 17757.666
$ws_ALC.Range("L86").Value = 68850
$ws_ALC.Range("M86").Value = -16634.666
$ws_ALC.Range("N86").Value = -71096
$ws_ALC.Range("H87").Value = 13334.892
$ws_ALC.Range("J87").Value = 13334.892
$ws_ALC.Range("L87").Value = 13334.892
$ws_ALC.Range("N87").Value = -15830.892
$ws_ALC.Range("H89").Value = 34788.445
$ws_ALC.Range("I89").Value = 17757.666
$ws_ALC.Range("J89").Value = 68850
$ws_ALC.Range("K89").Value = 88788.33
$ws_ALC.Range("L89").Value = 344250
$ws_ALC.Range("M89").Value = -83172.33
$ws_ALC.Range("N89").Value = -355482
$ws_ALC.Range("H90").Value = 13334.892
$ws_ALC.Range("J90").Value = 13334.892
$ws_ALC.Range("L90").Value = 40004.676
$ws_ALC.Range("N90").Value = -52484.676
$ws_ALC.Range("H112").Value = 1263.129
$ws_ALC.Range("J112").Value = 1306.32
$ws_ALC.Range("L112").Value = 3918.96
$ws_ALC.Range("N112").Value = -6134.96
$ws_ALC.Range("H113").Value = 6552.925
$ws_ALC.Range("I113").Value = 3085.9
$ws_ALC.Range("J113").Value = 10019.95
$ws_ALC.Range("K113").Value = 3085.9
$ws_ALC.Range("L113").Value = 10019.95
$ws_ALC.Range("M113").Value = 168.0999999999999
$ws_ALC.Range("N113").Value = -16527.95
$ws_ALC.Range("H129").Value = 1033.3617
$ws_ALC.Range("J129").Value = 1127.5641
$ws_ALC.Range("L129").Value = 3382.6923
$ws_ALC.Range("N129").Value = -13382.6923
$ws_ALC.Range("H132").Value = 38258.758
$ws_ALC.Range("I132").Value = 60405.777
$ws_ALC.Range("J132").Value = 2018.1818
$ws_ALC.Range("K132").Value = 181217.331
$ws_ALC.Range("L132").Value = 6054.5454
$ws_ALC.Range("M132").Value = -178687.331
$ws_ALC.Range("N132").Value = -11114.5454
$ws_ALC.Range("H138").Value = 3164.64
$ws_ALC.Range("I138").Value = 1974.0851
$ws_ALC.Range("J138").Value = 4220.415
$ws_ALC.Range("K138").Value = 5922.2553
$ws_ALC.Range("L138").Value = 12661.245
$ws_ALC.Range("M138").Value = -782.2552999999998
$ws_ALC.Range("N138").Value = -22941.245
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 6986.02
$ws_ARM.Range("I32").Value = 6072.0747
$ws_ARM.Range("J32").Value = 21304.5
$ws_ARM.Range("K32").Value = 6072.0747
$ws_ARM.Range("L32").Value = 21304.5
$ws_ARM.Range("M32").Value = -5785.0747
$ws_ARM.Range("N32").Value = -21878.5
$ws_ARM.Range("H110").Value = 50483.2
$ws_ARM.Range("I110").Value = 56036.89
$ws_ARM.Range("J110").Value = 500
$ws_ARM.Range("K110").Value = 56036.89
$ws_ARM.Range("L110").Value = 500
$ws_ARM.Range("M110").Value = -53991.89
$ws_ARM.Range("N110").Value = -4590
$ws_ARM.Range("H132").Value = 1854.1765
$ws_ARM.Range("I132").Value = 1367.4783
$ws_ARM.Range("J132").Value = 2871.818
$ws_ARM.Range("K132").Value = 4102.4349
$ws_ARM.Range("L132").Value = 8615.454000000002
$ws_ARM.Range("M132").Value = -1572.4349
$ws_ARM.Range("N132").Value = -13675.454
$ws_ARM.Range("H139").Value = 61806.3
$ws_ARM.Range("J139").Value = 61806.3
$ws_ARM.Range("L139").Value = 61806.3
$ws_ARM.Range("N139").Value = -72086.3
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 29752
$ws_BSM.Range("I20").Value = 38002.668
$ws_BSM.Range("J20").Value = 5000
$ws_BSM.Range("K20").Value = 38002.668
$ws_BSM.Range("L20").Value = 5000
$ws_BSM.Range("M20").Value = -37755.668
$ws_BSM.Range("N20").Value = -5494
$ws_BSM.Range("H134").Value = 1683.1
$ws_BSM.Range("I134").Value = 1272.6364
$ws_BSM.Range("K134").Value = 3817.9092
$ws_BSM.Range("M134").Value = -1282.9092
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 4186.156
$ws_CRP.Range("I31").Value = 2016.3103
$ws_CRP.Range("J31").Value = 8119
$ws_CRP.Range("K31").Value = 2016.3103
$ws_CRP.Range("L31").Value = 8119
$ws_CRP.Range("M31").Value = -1721.3103
$ws_CRP.Range("N31").Value = -8709
$ws_CRP.Range("H34").Value = 4186.156
$ws_CRP.Range("I34").Value = 2016.3103
$ws_CRP.Range("J34").Value = 8119
$ws_CRP.Range("K34").Value = 2016.3103
$ws_CRP.Range("L34").Value = 8119
$ws_CRP.Range("M34").Value = -1814.3103
$ws_CRP.Range("N34").Value = -8523
$ws_CRP.Range("H58").Value = 2602.1353
$ws_CRP.Range("I58").Value = 1684.5883
$ws_CRP.Range("J58").Value = 3382.05
$ws_CRP.Range("K58").Value = 1684.5883
$ws_CRP.Range("L58").Value = 3382.05
$ws_CRP.Range("M58").Value = -1481.5883
$ws_CRP.Range("N58").Value = -3788.05
$ws_CRP.Range("H134").Value = 2353.2593
$ws_CRP.Range("I134").Value = 2271.25
$ws_CRP.Range("J134").Value = 3009.3333
$ws_CRP.Range("K134").Value = 6813.75
$ws_CRP.Range("L134").Value = 9027.999899999999
$ws_CRP.Range("M134").Value = -4278.75
$ws_CRP.Range("N134").Value = -14097.9999
$ws_CRP.Range("H136").Value = 2602.1353
$ws_CRP.Range("I136").Value = 1684.5883
$ws_CRP.Range("J136").Value = 3382.05
$ws_CRP.Range("K136").Value = 5053.7649
$ws_CRP.Range("L136").Value = 10146.15
$ws_CRP.Range("M136").Value = -2503.7649
$ws_CRP.Range("N136").Value = -15246.15
$ws_CRP.Range("H140").Value = 76322.42999999999
$ws_CRP.Range("J140").Value = 76322.42999999999
$ws_CRP.Range("L140").Value = 76322.42999999999
$ws_CRP.Range("N140").Value = -86682.42999999999
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H107").Value = 378588.28
$ws_CUL.Range("I107").Value = 353.58334
$ws_CUL.Range("J107").Value = 567705.6
$ws_CUL.Range("K107").Value = 1060.75002
$ws_CUL.Range("L107").Value = 1703116.8
$ws_CUL.Range("M107").Value = 859.2499800000001
$ws_CUL.Range("N107").Value = -1706956.8
$ws_CUL.Range("H122").Value = 1179.6316
$ws_CUL.Range("I122").Value = 612.58826
$ws_CUL.Range("J122").Value = 5999.5
$ws_CUL.Range("K122").Value = 5513.29434
$ws_CUL.Range("L122").Value = 53995.5
$ws_CUL.Range("M122").Value = -3063.29434
$ws_CUL.Range("N122").Value = -58895.5
$ws_CUL.Range("H131").Value = 772.51
$ws_CUL.Range("I131").Value = 344.2857
$ws_CUL.Range("J131").Value = 842.22095
$ws_CUL.Range("K131").Value = 1032.8571
$ws_CUL.Range("L131").Value = 2526.66285
$ws_CUL.Range("M131").Value = 4007.1429
$ws_CUL.Range("N131").Value = -12606.66285
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H2").Value = 1298.4375
$ws_GSM.Range("I2").Value = 1303.125
$ws_GSM.Range("J2").Value = 1293.75
$ws_GSM.Range("K2").Value = 1303.125
$ws_GSM.Range("L2").Value = 1293.75
$ws_GSM.Range("M2").Value = -1190.125
$ws_GSM.Range("N2").Value = -1519.75
$ws_GSM.Range("H113").Value = 1344.8948
$ws_GSM.Range("I113").Value = 1131.5
$ws_GSM.Range("K113").Value = 1131.5
$ws_GSM.Range("M113").Value = 1038.5
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H127").Value = 41492.5
$ws_LTW.Range("J127").Value = 41492.5
$ws_LTW.Range("L127").Value = 41492.5
$ws_LTW.Range("N127").Value = -51412.5
$ws_LTW.Range("H136").Value = 3064.7856
$ws_LTW.Range("I136").Value = 3131.5293
$ws_LTW.Range("J136").Value = 2961.6365
$ws_LTW.Range("K136").Value = 9394.5879
$ws_LTW.Range("L136").Value = 8884.9095
$ws_LTW.Range("M136").Value = -6844.5879
$ws_LTW.Range("N136").Value = -13984.9095
$ws_LTW.Range("H140").Value = 56366.715
$ws_LTW.Range("J140").Value = 56366.715
$ws_LTW.Range("L140").Value = 56366.715
$ws_LTW.Range("N140").Value = -66726.715
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H42").Value = 34875
$ws_WVR.Range("J42").Value = 34875
$ws_WVR.Range("L42").Value = 34875
$ws_WVR.Range("N42").Value = -35631
$ws_WVR.Range("H121").Value = 23500
$ws_WVR.Range("J121").Value = 23500
$ws_WVR.Range("L121").Value = 23500
$ws_WVR.Range("N121").Value = -26994
$ws_WVR.Range("H132").Value = 2269.225
$ws_WVR.Range("I132").Value = 1695.5518
$ws_WVR.Range("J132").Value = 3781.6365
$ws_WVR.Range("K132").Value = 5086.6554
$ws_WVR.Range("L132").Value = 11344.9095
$ws_WVR.Range("M132").Value = -2556.6554
$ws_WVR.Range("N132").Value = -16404.9095
